# geração de análises seriais
# Re-sort the tail (tied/zero-value) rows of the ranking sheets so that the
# "uf" labels in column A appear in the newly generated order. Only the
# column A text values move; the underlying numeric values in column B are
# left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "tot-arrecad": rows 24-26 ---
$ws = $wb.Worksheets.Item("tot-arrecad")
$ws.Range("A24").Value = "RN"
$ws.Range("A25").Value = "AL"
$ws.Range("A26").Value = "PI"

# --- Sheet "avg-arrecad": rows 22-27 ---
$ws = $wb.Worksheets.Item("avg-arrecad")
$ws.Range("A22").Value = "PI"
$ws.Range("A23").Value = "MT"
$ws.Range("A24").Value = "RN"
$ws.Range("A25").Value = "RO"
$ws.Range("A26").Value = "AP"
$ws.Range("A27").Value = "TO"

# --- Sheet "max-arrecad": rows 21-27 ---
$ws = $wb.Worksheets.Item("max-arrecad")
$ws.Range("A21").Value = "AL"
$ws.Range("A22").Value = "MT"
$ws.Range("A23").Value = "RO"
$ws.Range("A24").Value = "AP"
$ws.Range("A25").Value = "TO"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"

# --- Sheet "tx-sucesso": rows 21-26 ---
$ws = $wb.Worksheets.Item("tx-sucesso")
$ws.Range("A21").Value = "RN"
$ws.Range("A22").Value = "AL"
$ws.Range("A23").Value = "RO"
$ws.Range("A24").Value = "AP"
$ws.Range("A25").Value = "TO"
$ws.Range("A26").Value = "MT"
